# Delete record 11431839 from the workbook.
# The record appears as a full data row in two worksheets:
#   - "Kayitlar"    (master list) -> row 1510
#   - "Merkez İlçe" (district-filtered view) -> row 971
# Deleting the entire row shifts all following rows up by one,
# which matches the target diff (no other edits are required).

$wb = $excel.ActiveWorkbook

$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$wsKayitlar.Rows.Item(1510).Delete()

$wsMerkez = $wb.Worksheets.Item("Merkez İlçe")
$wsMerkez.Rows.Item(971).Delete()
